# Fill in the "Responsible persons (2 Initials)" column (C) for the rows
# that were missing it on the "P2 worksheet" sheet, per the commit
# "Write output file, implementing DFA Tables, tests".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("P2 worksheet")

# Row 10 - "Structure of the program"
$ws.Range("C10").Value = "MA, AS"

# Row 13 - "The automaton (or several)" (Code and implementation section)
$ws.Range("C13").Value = "AS"

# Row 14 - "Read in the input file"
$ws.Range("C14").Value = "MA"

# Row 16 - "Create output file /release"
$ws.Range("C16").Value = "AS"

# Row 20 - "Comments in the code"
$ws.Range("C20").Value = "ALL"

# Row 21 - "Structure and readibility of the code"
$ws.Range("C21").Value = "ALL"

# Leave the cursor where the author left it when they saved the file.
$ws.Range("D15").Select()
